$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 9).Value = 1808
$ws.Cells.Item(3, 8).Value = 8349
$ws.Cells.Item(3, 9).Value = 1909
$ws.Cells.Item(4, 3).Value = 1804
$ws.Cells.Item(4, 5).Value = 1960
$ws.Cells.Item(4, 7).Value = 1424
$ws.Cells.Item(4, 8).Value = 1658
$ws.Cells.Item(4, 9).Value = 478
$ws.Cells.Item(6, 9).Value = 2347
$ws.Cells.Item(7, 3).Value = 28347
$ws.Cells.Item(7, 5).Value = 25964
$ws.Cells.Item(7, 7).Value = 24643
$ws.Cells.Item(7, 8).Value = 25970
$ws.Cells.Item(7, 9).Value = 6708

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 9).Value = 67
$ws.Cells.Item(7, 9).Value = 236
$ws.Cells.Item(8, 9).Value = 416
$ws.Cells.Item(11, 9).Value = 119
$ws.Cells.Item(13, 9).Value = 7
$ws.Cells.Item(14, 9).Value = 35
$ws.Cells.Item(15, 9).Value = 85
$ws.Cells.Item(18, 9).Value = 52
$ws.Cells.Item(19, 9).Value = 195
$ws.Cells.Item(20, 7).Value = 619
$ws.Cells.Item(20, 9).Value = 182
$ws.Cells.Item(23, 9).Value = 59
$ws.Cells.Item(27, 9).Value = 60
$ws.Cells.Item(29, 9).Value = 427
$ws.Cells.Item(31, 9).Value = 67
$ws.Cells.Item(33, 9).Value = 307
$ws.Cells.Item(34, 9).Value = 29
$ws.Cells.Item(37, 9).Value = 213
$ws.Cells.Item(42, 9).Value = 223
$ws.Cells.Item(43, 9).Value = 60
$ws.Cells.Item(47, 9).Value = 49
$ws.Cells.Item(48, 9).Value = 64
$ws.Cells.Item(49, 9).Value = 39
$ws.Cells.Item(52, 5).Value = 552
$ws.Cells.Item(52, 9).Value = 137
$ws.Cells.Item(55, 9).Value = 77
$ws.Cells.Item(56, 9).Value = 11
$ws.Cells.Item(59, 9).Value = 13
$ws.Cells.Item(60, 9).Value = 37
$ws.Cells.Item(61, 9).Value = 6
$ws.Cells.Item(63, 3).Value = 238
$ws.Cells.Item(63, 8).Value = 186
$ws.Cells.Item(63, 9).Value = 30
$ws.Cells.Item(65, 9).Value = 156
$ws.Cells.Item(67, 9).Value = 260
$ws.Cells.Item(72, 9).Value = 24
$ws.Cells.Item(73, 9).Value = 63
$ws.Cells.Item(75, 9).Value = 27
$ws.Cells.Item(76, 9).Value = 109
$ws.Cells.Item(78, 9).Value = 88
$ws.Cells.Item(79, 9).Value = 170
$ws.Cells.Item(81, 9).Value = 7
$ws.Cells.Item(85, 9).Value = 318
$ws.Cells.Item(87, 9).Value = 7
$ws.Cells.Item(89, 9).Value = 67
$ws.Cells.Item(90, 9).Value = 79
$ws.Cells.Item(96, 9).Value = 89
$ws.Cells.Item(100, 9).Value = 7
$ws.Cells.Item(101, 3).Value = 28347
$ws.Cells.Item(101, 5).Value = 25964
$ws.Cells.Item(101, 7).Value = 24643
$ws.Cells.Item(101, 8).Value = 25970
$ws.Cells.Item(101, 9).Value = 6708

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(3, 9).Value = 123
$ws.Cells.Item(7, 9).Value = 318

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 9).Value = 35
$ws.Cells.Item(3, 9).Value = 55
$ws.Cells.Item(4, 5).Value = 38
$ws.Cells.Item(4, 9).Value = 19
$ws.Cells.Item(7, 5).Value = 552
$ws.Cells.Item(7, 9).Value = 137

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(3, 9).Value = 23
$ws.Cells.Item(7, 9).Value = 119

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(4, 9).Value = 24
$ws.Cells.Item(6, 9).Value = 143
$ws.Cells.Item(7, 9).Value = 416

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 9).Value = 79
$ws.Cells.Item(6, 9).Value = 61
$ws.Cells.Item(7, 9).Value = 236

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(4, 9).Value = 9
$ws.Cells.Item(7, 9).Value = 67

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(6, 9).Value = 36
$ws.Cells.Item(7, 9).Value = 89

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Cells.Item(2, 9).Value = 10
$ws.Cells.Item(7, 9).Value = 35

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 9).Value = 68
$ws.Cells.Item(3, 9).Value = 63
$ws.Cells.Item(7, 9).Value = 213

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 9).Value = 60
$ws.Cells.Item(3, 9).Value = 89
$ws.Cells.Item(6, 9).Value = 91
$ws.Cells.Item(7, 9).Value = 260

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 9).Value = 19
$ws.Cells.Item(3, 9).Value = 20
$ws.Cells.Item(6, 9).Value = 24
$ws.Cells.Item(7, 9).Value = 67

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 9).Value = 40
$ws.Cells.Item(3, 9).Value = 43
$ws.Cells.Item(7, 9).Value = 156

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 9).Value = 71
$ws.Cells.Item(3, 9).Value = 106
$ws.Cells.Item(6, 9).Value = 108
$ws.Cells.Item(7, 9).Value = 307

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(6, 9).Value = 20
$ws.Cells.Item(7, 9).Value = 39

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 9).Value = 129
$ws.Cells.Item(3, 9).Value = 138
$ws.Cells.Item(6, 9).Value = 138
$ws.Cells.Item(7, 9).Value = 427

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(3, 9).Value = 47
$ws.Cells.Item(6, 9).Value = 55
$ws.Cells.Item(7, 9).Value = 195

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(4, 9).Value = 5
$ws.Cells.Item(6, 9).Value = 34
$ws.Cells.Item(7, 9).Value = 64

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(6, 9).Value = 46
$ws.Cells.Item(7, 9).Value = 109

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 9).Value = 60
$ws.Cells.Item(7, 9).Value = 223

$ws = $wb.Worksheets.Item('Boystown')
$ws.Cells.Item(3, 9).Value = 2
$ws.Cells.Item(6, 9).Value = 7

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(4, 9).Value = 17
$ws.Cells.Item(6, 9).Value = 32
$ws.Cells.Item(7, 9).Value = 88

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(3, 9).Value = 18
$ws.Cells.Item(7, 9).Value = 77

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(2, 9).Value = 17
$ws.Cells.Item(7, 9).Value = 59

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(3, 9).Value = 47
$ws.Cells.Item(7, 9).Value = 170

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(3, 9).Value = 55
$ws.Cells.Item(4, 7).Value = 25
$ws.Cells.Item(7, 7).Value = 619
$ws.Cells.Item(7, 9).Value = 182

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(6, 9).Value = 27
$ws.Cells.Item(7, 9).Value = 52

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Cells.Item(5, 9).Value = 4
$ws.Cells.Item(6, 9).Value = 7

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(2, 9).Value = 11
$ws.Cells.Item(7, 9).Value = 29

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(2, 9).Value = 7
$ws.Cells.Item(7, 9).Value = 49

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(2, 9).Value = 25
$ws.Cells.Item(3, 9).Value = 17
$ws.Cells.Item(6, 9).Value = 35
$ws.Cells.Item(7, 9).Value = 85

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(3, 9).Value = 20
$ws.Cells.Item(6, 9).Value = 16
$ws.Cells.Item(7, 9).Value = 63

$ws = $wb.Worksheets.Item('Montclare')
$ws.Cells.Item(3, 9).Value = 2
$ws.Cells.Item(7, 9).Value = 13

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(3, 9).Value = 24
$ws.Cells.Item(7, 9).Value = 67

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(3, 9).Value = 10
$ws.Cells.Item(6, 9).Value = 27
$ws.Cells.Item(7, 9).Value = 60

$ws = $wb.Worksheets.Item('Pullman')
$ws.Cells.Item(2, 9).Value = 11
$ws.Cells.Item(7, 9).Value = 27

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(2, 9).Value = 25
$ws.Cells.Item(7, 9).Value = 79

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(2, 9).Value = 5
$ws.Cells.Item(7, 9).Value = 37

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(3, 9).Value = 11
$ws.Cells.Item(7, 9).Value = 60

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(7, 9).Value = 24

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Cells.Item(6, 9).Value = 8
$ws.Cells.Item(7, 9).Value = 11

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(3, 9).Value = 4
$ws.Cells.Item(7, 9).Value = 7

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Cells.Item(6, 9).Value = 3
$ws.Cells.Item(7, 9).Value = 6

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Cells.Item(3, 9).Value = 2
$ws.Cells.Item(6, 9).Value = 7
